$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for numeric-looking price cells so they are not
# coerced into Excel numbers (the source data stores these as text strings).
$textCells = @("D5", "D8", "D9", "D12", "D18", "D19", "D20", "D26", "D31", "D34", "D38", "D39", "D40", "D44", "D45", "D46", "D47", "D48", "D51", "D42", "D43")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.594.37"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "1.597.04"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").Value = "212.08"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("D8").Value = "26.83"
$ws.Range("E8").Value = "  +4.94%  "

$ws.Range("D9").Value = "43.75"
$ws.Range("E9").Value = "  -4.11%  "

$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").Value = "0.0907"

$ws.Range("D13").Value = "1.824.92"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "1.591.44"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").Value = "29.592.59"
$ws.Range("E15").Value = "  +1.65%  "

$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "63.75"
$ws.Range("E18").Value = "  +2.01%  "

$ws.Range("D19").Value = "241.77"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  +2.22%  "

$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("D26").Value = "154.55"
$ws.Range("E26").Value = "  +1.00%  "

$ws.Range("E27").Value = "  +1.69%  "

$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  +1.08%  "

$ws.Range("E30").Value = "  +0.41%  "

$ws.Range("D31").Value = "0.0476"
$ws.Range("E31").Value = "  +2.61%  "

$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  +3.08%  "

$ws.Range("D35").Value = "1.431.27"
$ws.Range("E35").Value = "  +0.70%  "

$ws.Range("E36").Value = "  +1.95%  "

$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("D38").Value = "2.86"
$ws.Range("E38").Value = "  +2.52%  "

$ws.Range("D39").Value = "2.29"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "0.0165"
$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("E41").Value = "  +2.87%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "1.96"
$ws.Range("E42").Value = "  +1.06%  "

$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "54.56"
$ws.Range("E43").Value = "  +2.38%  "

$ws.Range("D44").Value = "0.0491"
$ws.Range("E44").Value = "  +6.65%  "

$ws.Range("D45").Value = "0.800"
$ws.Range("E45").Value = "  +1.57%  "

$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("D47").Value = "0.987"
$ws.Range("E47").Value = "  +16.10%  "

$ws.Range("D48").Value = "65.59"
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").Value = "1.736.92"
$ws.Range("E50").Value = "  +0.95%  "

$ws.Range("D51").Value = "85.92"
$ws.Range("E51").Value = "  +0.39%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}